$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column E (Patch Note Commentary) with bulleted, merged text ---
# Rows 2-4: bullets for the "Prompts UI..." note, row 2 also gets "First Patch."
$ws.Range("E2").Value = "• First Patch. `n• Prompts UI is messed up (ON DOWNPATCHED VERSION, it was fine back in the day) (Fixed February 2021)"
$ws.Range("E3").Value = "• Prompts UI is messed up (ON DOWNPATCHED VERSION, it was fine back in the day) (Fixed February 2021)"
$ws.Range("E4").Value = "• Prompts UI is messed up (ON DOWNPATCHED VERSION, it was fine back in the day) (Fixed February 2021)"

$ws.Range("E5").Value = "• Removed Druzhina wallbang`n• Removed Train OOB`n• Removed Berlin Manhole Exit`n• Removed a Chongqing Vault Skip`n• Weird downpatch prompt issue fixed"

$ws.Range("E7").Value = "• Check `"April 6 - PC Only`" in official patch notes. Fixed a frame rate issue."

$ws.Range("E8").Value = "• Emetic Patch (sick NPC throw up elsewhere), Legal Sniper in Mendoza (Gaucho Start) patched in Mendoza. Also some lighting fixes."

$ws.Range("E10").Value = "• Check `"June 21 - PC Only`" in official patch notes. Fixed an issue with Epic Overlay."

$ws.Range("E11").Value = "• Added Dartmoor game show"

$ws.Range("E12").Value = "• Check `"August 31 - PC Only`" in official patch notes. Fixed an issue where players experienced FPS drops/ stuttering after the 3.50 (July) patch."

$ws.Range("E13").Value = "• Added back shoulder swap. Tweaked sniper slowdown."

$ws.Range("E14").Value = "• Check `"October 1 - PC Only`" in official patch notes. Fixed `"No “Woosh” sound on PC. "

$ws.Range("E16").Value = "• Year 2 Patch. Removed Dartmoor SA 37 chandy wallbang and messed up sliding doors`n• First patch available to Steam users.`n• First patch Molotov seems to be available."

$ws.Range("E17").Value = "• Year 2 variant."

$ws.Range("E18").Value = "• Year 2 variant. `n• No extra cameras in Master difficulty Dubai yet. `n• May 19th according to featured contract thumbnails: https://www.hitmanforum.com/t/year-2-may-official-community-rubber-duck-featured-contracts/14702"

$ws.Range("E19").Value = "• S3 Master Cameras Patch`n• Added invisible floor to Mendoza sliding roof`n• Most recent ET Arcade dates to 2022-Jul-14"
$ws.Range("E20").Value = "• S3 Master Cameras Patch`n• Added invisible floor to Mendoza sliding roof`n• Most recent ET Arcade dates to 2022-Jul-14"

$ws.Range("E21").Value = "• Added Molotov. `n• Technically not the only patch with Molly playable"

$ws.Range("E22").Value = "• Removed molly accident / wallbang`n• Changed NY frisk`n• Changed Dartmoor ledge drop`n• Removed RFID exploit"

# --- Clear out the now-unused "Other Notes" column (F) data rows ---
$ws.Range("F2:F22").Clear()

# --- Row heights (content grew due to bullets / merged notes) ---
$ws.Rows.Item(2).RowHeight = 86.4
$ws.Rows.Item(3).RowHeight = 72
$ws.Rows.Item(4).RowHeight = 72
$ws.Rows.Item(5).RowHeight = 129.6
$ws.Rows.Item(16).RowHeight = 115.2
$ws.Rows.Item(18).RowHeight = 172.8
$ws.Rows.Item(19).RowHeight = 86.4
$ws.Rows.Item(20).RowHeight = 86.4
$ws.Rows.Item(22).RowHeight = 86.4

# --- Sheet view: scrolled down, new selection ---
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("F21").Select()
